$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update title text (A1) with new time
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 17:05"

# 2. Row 4 - Estados Unidos
$ws.Range("B4").Value = 1839698
$ws.Range("C4").Value = 2528
$ws.Range("D4").Value = 599925
$ws.Range("E4").Value = 1133511
$ws.Range("G4").Value = 67
$ws.Range("H4").Value = 106262

# 3. Row 10 - India
$ws.Range("B10").Value = 194504
$ws.Range("C10").Value = 3895
$ws.Range("D10").Value = 93343
$ws.Range("E10").Value = 95713
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 5448

# 4. Row 45 - Argentina
$ws.Range("D45").Value = 5521
$ws.Range("E45").Value = 10789
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 541

# 5. Row 62 - Moldavia
$ws.Range("B62").Value = 8360
$ws.Range("C62").Value = 109
$ws.Range("E62").Value = 3433
$ws.Range("G62").Value = 10
$ws.Range("H62").Value = 305

# 6. Row 68 - Irak
$ws.Range("B68").Value = 6868
$ws.Range("C68").Value = 429
$ws.Range("D68").Value = 3275
$ws.Range("E68").Value = 3378
$ws.Range("G68").Value = 10
$ws.Range("H68").Value = 215

# 7. Row 70 - Azerbaiyan
$ws.Range("B70").Value = 5662
$ws.Range("C70").Value = 168
$ws.Range("D70").Value = 3508
$ws.Range("E70").Value = 2086
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 68

# 8. Rows 93/94 - swap Kenia and Somalia, with Kenia getting new updated numbers
$ws.Range("A93").Value = "Kenia"
$ws.Range("B93").Value = 2021
$ws.Range("C93").Value = 59
$ws.Range("D93").Value = 478
$ws.Range("E93").Value = 1479
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 64

$ws.Range("A94").Value = "Somalia"
$ws.Range("B94").Value = 1976
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 348
$ws.Range("E94").Value = 1550
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 78
